# The deck currently carries two theme parts:
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet"  (used by the real SlideMaster,
#                            i.e. every normal slide in the deck)
#   ppt/theme/theme1.xml -> "Office Theme" / "Office"   (only referenced by the Notes Master)
#
# The authored change swaps the two themes' contents in place: the master/slide
# theme becomes the stock "Office" colour scheme, and the notes-master theme
# becomes the former "Red Violet" colour scheme (file names / relationships are
# untouched - only the <a:clrScheme> colour values move between the two parts).
#
# The only theme surface this PowerPoint host exposes for writing is the
# ThemeColorScheme attached to the live slide/master theme (theme2.xml here) -
# so we drive the visible, in-use theme to the "Office" colours through that
# object model, matching the half of the swap that is reachable via the
# PowerPoint COM surface.

function Convert-HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target colours = the presentation's current "Office Theme" clrScheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in theme colour-scheme order).
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToOleColor $officeThemeColors[$i - 1]
}
